$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 4 ("Couple and 2 children household ...")
# so the new household type is placed right after "Couple household".
$ws.Rows.Item(4).Insert()

# Fill in the new row's data: ID=3, name="Couple and 1 children household (both adults are fully-employeed)"
$ws.Cells.Item(4, 1).Value = 3
$ws.Cells.Item(4, 2).Value = "Couple and 1 children household (both adults are fully-employeed)"

# Renumber the IDs of the rows that shifted down.
$ws.Cells.Item(5, 1).Value = 4
$ws.Cells.Item(6, 1).Value = 5

# Match the author's final selection state.
$ws.Range("B8").Select() | Out-Null
